$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Add the new backlog item as a new row at the bottom of the task list
$ws.Range("A15").Value = "reporte ot por sector , agregar filtro por estados agregar columna de cantidad y NC"
$ws.Range("B15").Value = "no comenzado"

# Update the selected cell to reflect where the user left off after adding the row
$ws.Range("A16").Select()
